$d = $word.ActiveDocument

# Locate the "PlainText = EditText => PlainText Input" paragraph and
# underline the whole paragraph (all runs plus the paragraph mark).
$rng = $d.Content
$found = $rng.Find.Execute("PlainText = EditText => PlainText Input", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $rng.Paragraphs(1)
    $para.Range.Underline = 1
}
